$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1939.8572
$ws.Range("J97").Value = 1939.8572
$ws.Range("L97").Value = 5819.571599999999
$ws.Range("N97").Value = -6811.571599999999
$ws.Range("H111").Value = 2104.5
$ws.Range("J111").Value = 2905.5
$ws.Range("L111").Value = 8716.5
$ws.Range("N111").Value = -14850.5
$ws.Range("H112").Value = 2164.238
$ws.Range("I112").Value = 3200.3333
$ws.Range("J112").Value = 1991.5555
$ws.Range("K112").Value = 9600.999899999999
$ws.Range("L112").Value = 5974.666499999999
$ws.Range("M112").Value = -8492.999899999999
$ws.Range("N112").Value = -8190.666499999999
$ws.Range("H121").Value = 722.94446
$ws.Range("J121").Value = 722.94446
$ws.Range("L121").Value = 2168.83338
$ws.Range("N121").Value = -5662.83338
$ws.Range("H132").Value = 2175.182
$ws.Range("I132").Value = 2251.1614
$ws.Range("J132").Value = 997.5
$ws.Range("K132").Value = 6753.4842
$ws.Range("L132").Value = 2992.5
$ws.Range("M132").Value = -4223.4842
$ws.Range("N132").Value = -8052.5
$ws.Range("H137").Value = 6648.15
$ws.Range("I137").Value = 1555.3334
$ws.Range("J137").Value = 14287.375
$ws.Range("K137").Value = 4666.0002
$ws.Range("L137").Value = 42862.125
$ws.Range("M137").Value = -2116.0002
$ws.Range("N137").Value = -47962.125
$ws.Range("H138").Value = 2991.1555
$ws.Range("J138").Value = 3376.1714
$ws.Range("L138").Value = 10128.5142
$ws.Range("N138").Value = -20408.5142
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1507.9166
$ws.Range("I63").Value = 1507.9166
$ws.Range("K63").Value = 1507.9166
$ws.Range("M63").Value = -821.9166
$ws.Range("H66").Value = 1507.9166
$ws.Range("I66").Value = 1507.9166
$ws.Range("K66").Value = 7539.583000000001
$ws.Range("M66").Value = -4107.583000000001
$ws.Range("H74").Value = 2798.4
$ws.Range("I74").Value = 2599.8
$ws.Range("K74").Value = 2599.8
$ws.Range("M74").Value = -1725.8
$ws.Range("H77").Value = 2798.4
$ws.Range("I77").Value = 2599.8
$ws.Range("K77").Value = 12999
$ws.Range("M77").Value = -8631
$ws.Range("H97").Value = 49812.445
$ws.Range("J97").Value = 174686
$ws.Range("L97").Value = 174686
$ws.Range("N97").Value = -175678
$ws.Range("H132").Value = 3373.2727
$ws.Range("I132").Value = 3495.6833
$ws.Range("J132").Value = 2149.1667
$ws.Range("K132").Value = 10487.0499
$ws.Range("L132").Value = 6447.500100000001
$ws.Range("M132").Value = -7957.0499
$ws.Range("N132").Value = -11507.5001
$ws.Range("H138").Value = 88998.664
$ws.Range("J138").Value = 88998.664
$ws.Range("L138").Value = 88998.664
$ws.Range("N138").Value = -99278.664
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7164.769
$ws.Range("I107").Value = 6611.467
$ws.Range("J107").Value = 7919.273
$ws.Range("K107").Value = 6611.467
$ws.Range("L107").Value = 7919.273
$ws.Range("M107").Value = -4691.467
$ws.Range("N107").Value = -11759.273
$ws.Range("H134").Value = 3543.6775
$ws.Range("I134").Value = 2862.8572
$ws.Range("K134").Value = 8588.571599999999
$ws.Range("M134").Value = -6053.571599999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3235.2
$ws.Range("I31").Value = 3150.2778
$ws.Range("J31").Value = 3999.5
$ws.Range("K31").Value = 3150.2778
$ws.Range("L31").Value = 3999.5
$ws.Range("M31").Value = -2855.2778
$ws.Range("N31").Value = -4589.5
$ws.Range("H34").Value = 3235.2
$ws.Range("I34").Value = 3150.2778
$ws.Range("J34").Value = 3999.5
$ws.Range("K34").Value = 3150.2778
$ws.Range("L34").Value = 3999.5
$ws.Range("M34").Value = -2948.2778
$ws.Range("N34").Value = -4403.5
$ws.Range("H51").Value = 13332.667
$ws.Range("J51").Value = 13332.667
$ws.Range("L51").Value = 13332.667
$ws.Range("N51").Value = -14804.667
$ws.Range("H58").Value = 1452.5454
$ws.Range("I58").Value = 1342.9
$ws.Range("K58").Value = 1342.9
$ws.Range("M58").Value = -1139.9
$ws.Range("H61").Value = 13332.667
$ws.Range("J61").Value = 13332.667
$ws.Range("L61").Value = 13332.667
$ws.Range("N61").Value = -14028.667
$ws.Range("H107").Value = 1083.8125
$ws.Range("I107").Value = 807.6667
$ws.Range("J107").Value = 1912.25
$ws.Range("K107").Value = 807.6667
$ws.Range("L107").Value = 1912.25
$ws.Range("M107").Value = 1112.3333
$ws.Range("N107").Value = -5752.25
$ws.Range("H132").Value = 3133.7144
$ws.Range("I132").Value = 3098.36
$ws.Range("J132").Value = 3428.3333
$ws.Range("K132").Value = 9295.08
$ws.Range("L132").Value = 10284.9999
$ws.Range("M132").Value = -6765.08
$ws.Range("N132").Value = -15344.9999
$ws.Range("H134").Value = 2543.1924
$ws.Range("I134").Value = 2297
$ws.Range("K134").Value = 6891
$ws.Range("M134").Value = -4356
$ws.Range("H136").Value = 1452.5454
$ws.Range("I136").Value = 1342.9
$ws.Range("K136").Value = 4028.7
$ws.Range("M136").Value = -1478.7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5296496.5
$ws.Range("I4").Value = 1385554
$ws.Range("K4").Value = 4156662
$ws.Range("M4").Value = -4156550
$ws.Range("H11").Value = 500458.12
$ws.Range("I11").Value = 53099
$ws.Range("J11").Value = 3333732.8
$ws.Range("K11").Value = 159297
$ws.Range("L11").Value = 10001198.4
$ws.Range("M11").Value = -159157
$ws.Range("N11").Value = -10001478.4
$ws.Range("H14").Value = 611.25
$ws.Range("I14").Value = 611.25
$ws.Range("K14").Value = 1833.75
$ws.Range("M14").Value = -1660.75
$ws.Range("H92").Value = 1859.6666
$ws.Range("I92").Value = 700
$ws.Range("J92").Value = 2439.5
$ws.Range("K92").Value = 2100
$ws.Range("L92").Value = 7318.5
$ws.Range("M92").Value = -852
$ws.Range("N92").Value = -9814.5
$ws.Range("H113").Value = 357.7143
$ws.Range("I113").Value = 175
$ws.Range("J113").Value = 494.75
$ws.Range("K113").Value = 525
$ws.Range("L113").Value = 1484.25
$ws.Range("M113").Value = 1645
$ws.Range("N113").Value = -5824.25
$ws.Range("H122").Value = 408
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H127").Value = 4716.6665
$ws.Range("J127").Value = 4716.6665
$ws.Range("L127").Value = 14149.9995
$ws.Range("N127").Value = -24069.9995
$ws.Range("H138").Value = 3286.8125
$ws.Range("I138").Value = 1289.909
$ws.Range("K138").Value = 3869.727
$ws.Range("M138").Value = 1270.273
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 347.94736
$ws.Range("I2").Value = 307.46667
$ws.Range("J2").Value = 499.75
$ws.Range("K2").Value = 307.46667
$ws.Range("L2").Value = 499.75
$ws.Range("M2").Value = -194.46667
$ws.Range("N2").Value = -725.75
$ws.Range("H126").Value = 50062
$ws.Range("I126").Value = 4458.4
$ws.Range("J126").Value = 164071
$ws.Range("K126").Value = 13375.2
$ws.Range("L126").Value = 492213
$ws.Range("M126").Value = -10905.2
$ws.Range("N126").Value = -497153
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2614625
$ws.Range("J20").Value = 3419500
$ws.Range("L20").Value = 3419500
$ws.Range("N20").Value = -3419952
$ws.Range("H61").Value = 22248042
$ws.Range("I61").Value = 30306040
$ws.Range("K61").Value = 30306040
$ws.Range("M61").Value = -30305838
$ws.Range("H113").Value = 22248042
$ws.Range("I113").Value = 30306040
$ws.Range("K113").Value = 30306040
$ws.Range("M113").Value = -30303870
$ws.Range("H132").Value = 5730
$ws.Range("I132").Value = 4995
$ws.Range("K132").Value = 14985
$ws.Range("M132").Value = -12455
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H107").Value = 33378410
$ws.Range("I107").Value = 14847
$ws.Range("J107").Value = 62571530
$ws.Range("K107").Value = 44541
$ws.Range("L107").Value = 187714590
$ws.Range("M107").Value = -42621
$ws.Range("N107").Value = -187718430
$ws.Range("H113").Value = 1952.7
$ws.Range("I113").Value = 1517.8572
$ws.Range("K113").Value = 4553.571599999999
$ws.Range("M113").Value = -2383.571599999999
$ws.Range("H136").Value = 4038.4583
$ws.Range("I136").Value = 4458.7
$ws.Range("J136").Value = 1937.25
$ws.Range("K136").Value = 13376.1
$ws.Range("L136").Value = 5811.75
$ws.Range("M136").Value = -10826.1
$ws.Range("N136").Value = -10911.75
